# "completed all in reactJS assignment"
# Add the 28/08/2016 tracker row: Reply and React Router are now marked Completed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "28/08/2016"
$ws.Range("K9").Value = "Completed"
$ws.Range("L9").Value = "Completed"

# Keep the selection in sync with the new last row, like the saved workbook.
$ws.Range("M9").Select()
